$d = $word.ActiveDocument

function Find-LastMatch($doc, [string]$text) {
    $searchStart = 0
    $lastStart = -1
    $lastEnd = -1
    while ($true) {
        $r = $doc.Range($searchStart, $doc.Content.End)
        $found = $r.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if (-not $found) { break }
        $lastStart = $r.Start
        $lastEnd = $r.End
        $searchStart = $r.End
    }
    if ($lastStart -eq -1) { return $null }
    return $doc.Range($lastStart, $lastEnd)
}

# -----------------------------------------------------------------
# 1) Insert a new paragraph "Vi förväntar oss..." right after the
#    "Nedan presenteras fynd..." paragraph near the top of the doc.
# -----------------------------------------------------------------
$introText = "Nedan presenteras fynd av naturvårdsarter och fridlysta arter som gjorts i det avverkningsanmälda området, samt relevanta utdrag ur standarderna för FSC, Chain of Custody, Controlled Wood och PEFC. I BILAGA 1 finns artfakta om fridlysta arter."
$expectText = "Vi förväntar oss att ni återkommer med ett skriftligt svar på vårt klagomål och även beskriver vilka korrigerande åtgärder ni satt in för att rätta till identifierade brister i er efterlevnad av den svenska FSC standarden."

$introMatch = Find-LastMatch $d $introText
$introMatch.Collapse(0) | Out-Null
$introEndBeforeInsert = $introMatch.End
$introMatch.InsertParagraphAfter()

# The freshly inserted (empty) paragraph now sits right after "Nedan...".
# NB: $introMatch (being a collapsed range) does not track the edit made
# via InsertParagraphAfter, so its End still reports the pre-insert
# position; the new paragraph mark was placed right there, so the new
# paragraph itself begins one character further along.
$newParaStart = $introEndBeforeInsert + 1
$newPara = $d.Range($newParaStart, $newParaStart)
$newPara.Expand(4) | Out-Null                 # wdParagraph
$newParaText = $d.Range($newPara.Start, $newPara.End - 1)
$newParaText.Text = $expectText

# -----------------------------------------------------------------
# 2) Remove the old copy of that paragraph (and the two blank
#    paragraphs preceding it), located just before the page break,
#    right after the "...artskyddsförordningen" comment paragraph.
#    (There are two near-identical "Kommentar:" paragraphs with this
#    text in the document; we need the last/second one.)
# -----------------------------------------------------------------
$commentText = "I den avverkningsanmälda skogen har fridlysta arter sina livsmiljöer och växtplatser. Att skada de fridlysta arternas livsmiljöer, växtplatser eller ekologiska funktion är inte tillåtet enligt artskyddsförordningen"

$commentMatch = Find-LastMatch $d $commentText
$commentMatch.Expand(4) | Out-Null            # expand to the whole "Kommentar:" paragraph

# Walk forward three paragraphs: the two blank ones + the duplicate
# "Vi förväntar oss..." paragraph that should be removed.
$p1 = $d.Range($commentMatch.End, $commentMatch.End)
$p1.Expand(4) | Out-Null
$p2 = $d.Range($p1.End, $p1.End)
$p2.Expand(4) | Out-Null
$p3 = $d.Range($p2.End, $p2.End)
$p3.Expand(4) | Out-Null

$deleteRange = $d.Range($commentMatch.End, $p3.End)
$deleteRange.Delete()

# -----------------------------------------------------------------
# 3) Update the date in the header from 2023-11-13 to 2023-11-14.
# -----------------------------------------------------------------
foreach ($sec in $d.Sections) {
    foreach ($idx in 1, 2, 3) {
        $hdr = $sec.Headers.Item($idx)
        if ($hdr.Exists -and $hdr.Range.Text -like "*2023-11-13*") {
            $hdr.Range.Find.Execute("2023-11-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-14", 2) | Out-Null
        }
    }
}
